$d = $word.ActiveDocument

# Change 1: split the trailing period of the RLOC sentence into its own run.
$d.Content.Find.Execute(
    "e considering code smells are the CLOC, LOC and RLOC.", $true, $false, $false, $false, $false,
    $true, 1, $false, "e considering code smells are the CLOC, LOC and RLOC.", 2)

# Change 2: merge the "pretty excessive" run (and drop proofErr wrappers)
$d.Content.Find.Execute(
    "we have a method with 324 lines which is pretty excessive for a method.", $true, $false, $false, $false, $false,
    $true, 1, $false, "we have a method with 324 lines which is pretty excessive for a method.", 2)
